$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column D - copy formatting from the neighboring header (C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Ano"

# Fill the year value for each data row (2 through 10)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = 2024
}
